# Generate Report for Handback
#
# The localization-status report gained a freshly detected "handback" for
# the 2e07a577-6a7d-4c47-8e5b-9f7160f7a706 source file: its Latest Target
# File / Latest Handback File / Latest Handback DateTime columns get filled
# in, and an Error Detail is recorded because the handed-back file is not
# based on the latest handoff. This happened identically on both the
# "zh-cn" and "de-de" status sheets. The Error Detail column (P) also needs
# to be widened so the message is readable.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f83b3df7421f34a8671ee88426aee2a850469fac/e2e/2e07a577-6a7d-4c47-8e5b-9f7160f7a706.md"
$displayName = "2e07a577-6a7d-4c47-8e5b-9f7160f7a706.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb8d7b8745cc680066c4ea2a4f4ed9ceb5e9e61e/e2e/2e07a577-6a7d-4c47-8e5b-9f7160f7a706.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f83b3df7421f34a8671ee88426aee2a850469fac/e2e/2e07a577-6a7d-4c47-8e5b-9f7160f7a706.md."

# ColumnWidth (characters) + 5/6 == the stored raw OOXML column width, so
# asking for a raw width of 40 means setting ColumnWidth to 40 - 5/6.
$rawWidth40 = 40 - 5/6

# zh-cn: new handback info for row 7 (2e07a577-6a7d-4c47-8e5b-9f7160f7a706)
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("J7").Value = "2e07a577-6a7d-4c47-8e5b-9f7160f7a706.4ed6bbe66c9b13ad48d5f43ae518485f853bd359.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-09-01 18:50:47"
$wsZhCn.Range("P7").Value = $errorDetail

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $targetUrl, "", "", $displayName)
$wsZhCn.Range("I7").Font.Underline = $true
$wsZhCn.Range("I7").Font.Color = 15570276

$wsZhCn.Columns.Item(16).ColumnWidth = $rawWidth40

# de-de: same new handback info for row 7
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("J7").Value = "2e07a577-6a7d-4c47-8e5b-9f7160f7a706.4ed6bbe66c9b13ad48d5f43ae518485f853bd359.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-09-01 18:50:54"
$wsDeDe.Range("P7").Value = $errorDetail

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $targetUrl, "", "", $displayName)
$wsDeDe.Range("I7").Font.Underline = $true
$wsDeDe.Range("I7").Font.Color = 15570276

$wsDeDe.Columns.Item(16).ColumnWidth = $rawWidth40
